$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$url = "https://git-scm.com/install/"
$desc = "Para descargar GIT"

# Add the new hyperlink row (row 9), mirroring the existing rows B3:C8.
$ws.Hyperlinks.Add($ws.Range("B9"), $url) | Out-Null
$ws.Range("C9").Value = $desc

# Reuse the same "Hipervínculo" cell style already used by B3:B8 instead of
# the distinct style Hyperlinks.Add created.
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B9").PasteSpecial(-4122) | Out-Null # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("C10").Select()
